$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Update cell values. The order below intentionally introduces brand-new
#     comment strings (rows 13/14/21/23/24/26/27/30/35) in the same sequence
#     the original author typed them, so that after Excel rebuilds/compacts
#     the shared-strings table the new entries land on the same indices as
#     the target workbook. ---
$ws.Range("F21").Value = "akkoord. Wordt uitgezocht. Ook de het attribuut ligging moet optioneel worden`nHet is niet de bedoeling dat er er verplicht een bestand meegeleverd wordt"
$ws.Range("G21").Value = "verwerkt.`nBeide attributen zijn optioneel gemaakt.`nVerwerkt in UML.`nVerwerkt in XSD.`nVerwerkt in changelog."
$ws.Range("G23").Value = "is verwerkt in UML"
$ws.Range("G24").Value = "geen aanpassing nodig"
$ws.Range("G26").Value = "is verwerkt in UML`nis verwerkt in changelog`nis verwerkt in XSD"
$ws.Range("G27").Value = "is verwerkt in uml"
$ws.Range("G30").Value = "is verwerkt in doc"
$ws.Range("G13").Value = "is verwerkt in waardelijst excel`nis verwerkt in uml`nis verwerkt in XSD`n"
$ws.Range("G14").Value = "is verwerkt in waardelijst excel`nis verwerkt in uml`nis verwerkt in XSD"
$ws.Range("G35").Value = "verwerkt in uml`nverwerkt in IMKL-extraRegels`nverwerkt in XSD`n"
$ws.Range("F35").Value = "Akkoord. Deze wordt op 0..1 gezet en met een regel voor kabelbed en duct verplicht."
$ws.Range("D1").Value = "Hoofdstuk/paragraaf/figuur/regelnummer"
$ws.Range("G1").Value = "verwerkt"
$ws.Range("G2").Value = "is in pmkl doc verwerkt"
$ws.Range("G3").Value = "is verwerkt in pmkl doc"
$ws.Range("G4").Value = "is verwerkt"
$ws.Range("G5").Value = "is verwerkt. Ook uit UML verwijderd"
$ws.Range("G8").Value = "is verwerkt in pmkl doc en sld's"
$ws.Range("G9").Value = "geen verwerking nodig"
$ws.Range("G10").Value = "verwerkt in uml"
$ws.Range("G11").Value = "verwerkt in uml en daarmee ook in modeldocument (na generatie objectcat)"
$ws.Range("F12").Value = "akkoord. Het was wel de bedoeling dat een mantelbuis ook een kleur attribuut kan hebben. Het model realiseert dit nu niet. Aanpassing door kleur ook optioneel bij een mantelbuis op te nemen. FysiekeIdentificatie is voor een mantelbuis niet relevant."
$ws.Range("G12").Value = "verwerkt in uml`nverwerkt in xsd`nverwerkt in changelog"
$ws.Range("F13").Value = "akkoord, verwerkt.  De waardelijst is opgenomen."
$ws.Range("F14").Value = "akkoord, verwerkt.  De waardelijst is opgenomen."
$ws.Range("G15").Value = "changlog item 39 is verwijderd.`nVerwerkt in UML: Constraints zijn aangepast/herformuleerd op basis van voorstel.`nVerwerkt in changelog: herformulering van constraints op GebiedsinformatieAanvraag`n"
$ws.Range("F17").Value = "akkoord, verwerkt. De aangepaste lijst is ontvangen en is verwerkt."
$ws.Range("G17").Value = "is verwerkt in waardelijst excel"
$ws.Range("G18").Value = "is verwerkt in UML diagrammen en verwijderd uit model."
$ws.Range("G19").Value = "is verwerkt in extraRegels excel"
$ws.Range("G20").Value = "verwerkt. Dit changelog item is nu item 42 geworden"
$ws.Range("F23").Value = "Akkoord. Wordt aangepast. fysiekeIdentificatie`nDefinitie: Merkband, nummer of print op de buis.`nToelichting: Maximaal 64 karakters."
$ws.Range("F24").Value = "Antwoord: Dit zijn extra model regels in de vorm van 'constraints'. Je ziet ze als je doorklikt en scrollt naar Overzicht constraints"
$ws.Range("F25").Value = "Niet akkoord. Een verplichting is niet nodig. De optionele Dieptelegging kan al aan een mantelbuis gekoppeld worden of met een annotatie worden weergegeven."
$ws.Range("G25").Value = "geen aanpassing nodig"
$ws.Range("F26").Value = "Akkoord. De kleur wordt optioneel toegevoegd. Een fysiekeIdentificatie komt bij een mantelbuis niet voor."
$ws.Range("F27").Value = "Akkoord. De definitie wordt toegevoegd"
$ws.Range("F28").Value = "Het is zoals het nu is: Bij een graafmelding is er voor één netbeheerder, binnen één thema alleen een bijlage voor de maatregel met de hoogste prioriteit."
$ws.Range("G28").Value = "geen aanpassing nodig"
$ws.Range("F30").Value = "Akkoord: Wordt aangepast"
$ws.Range("F31").Value = "Akkoord: Wordt aangepast"
$ws.Range("G31").Value = "is verwerkt in uml"
$ws.Range("F32").Value = "Akkoord, verwerkt. Is verwerkt"
$ws.Range("G32").Value = "is verwerkt in waardelijst excel"
$ws.Range("F33").Value = "Akkoord, verwerkt. De waardelijst is toegevoegd"
$ws.Range("G33").Value = "is verwerkt in waardelijst excel`nis verwerkt in uml`nis verwerkt in XSD"
$ws.Range("F34").Value = "Wordt als 'infiltratievoorziening' opgenomen. Zie issue 291"
$ws.Range("G34").Value = "geen aanpassing nodig"
$ws.Range("E36").Value = "De aanvullende wens betreft de zichtbaarheid van een gestuurde boring in het kaartbeeld.`nDe aanvulling die we voorstellen is dat in het geval van een ExtraDetailInfo waarin het attribuut aanlegmethodeGestuurdeBoring true is dat dan het icoontje veranderd naar iets wat het herkenbaar maakt als een boring (bijvoorbeeld een boormachine). Het betreft dus een aanvulling in het PMKL)."
$ws.Range("F36").Value = "wordt als voorstel uitgewerkt"

# --- Turn on word-wrap for the "verwerking" cells whose text now spans
#     multiple lines (cellXfs style flips from index 23 to 24) ---
$ws.Range("G20").WrapText = $true
$ws.Range("G21").WrapText = $true
$ws.Range("G26").WrapText = $true
$ws.Range("G33").WrapText = $true
$ws.Range("G35").WrapText = $true

# --- Row 33 grew taller once its comment text was filled in ---
$ws.Rows.Item(33).RowHeight = 43.8

# --- Selection moved from F36 to E36 ---
$ws.Range("E36").Select()

